$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 4324.1113
$ws.Range("I58").Value = 150
$ws.Range("J58").Value = 5516.7144
$ws.Range("K58").Value = 450
$ws.Range("L58").Value = 16550.1432
$ws.Range("M58").Value = -300
$ws.Range("N58").Value = -16850.1432

$ws.Range("H64").Value = 3440.0845
$ws.Range("I64").Value = 3597.2156
$ws.Range("J64").Value = 3039.4
$ws.Range("K64").Value = 3597.2156
$ws.Range("L64").Value = 3039.4
$ws.Range("M64").Value = -3349.2156
$ws.Range("N64").Value = -3535.4

$ws.Range("H67").Value = 3440.0845
$ws.Range("I67").Value = 3597.2156
$ws.Range("J67").Value = 3039.4
$ws.Range("K67").Value = 3597.2156
$ws.Range("L67").Value = 3039.4
$ws.Range("M67").Value = -2739.2156
$ws.Range("N67").Value = -4755.4

$ws.Range("H80").Value = 479.33334
$ws.Range("J80").Value = 501.81818
$ws.Range("L80").Value = 1505.45454
$ws.Range("N80").Value = -3501.45454

$ws.Range("H83").Value = 479.33334
$ws.Range("J83").Value = 501.81818
$ws.Range("L83").Value = 4516.36362
$ws.Range("N83").Value = -14500.36362

$ws.Range("H113").Value = 5436638.5
$ws.Range("I113").Value = 1913.2858
$ws.Range("J113").Value = 13890655
$ws.Range("K113").Value = 1913.2858
$ws.Range("L113").Value = 13890655
$ws.Range("M113").Value = 1340.7142
$ws.Range("N113").Value = -13897163

$ws.Range("H116").Value = 8874.0625
$ws.Range("I116").Value = 12798.333
$ws.Range("J116").Value = 3828.5715
$ws.Range("K116").Value = 12798.333
$ws.Range("L116").Value = 3828.5715
$ws.Range("M116").Value = -9356.333000000001
$ws.Range("N116").Value = -10712.5715

$ws.Range("H129").Value = 1039.5541
$ws.Range("I129").Value = 712.6667
$ws.Range("J129").Value = 1102.8226
$ws.Range("K129").Value = 2138.0001
$ws.Range("L129").Value = 3308.4678
$ws.Range("M129").Value = 2861.9999
$ws.Range("N129").Value = -13308.4678

$ws.Range("H132").Value = 2558.65
$ws.Range("I132").Value = 2515
$ws.Range("K132").Value = 7545
$ws.Range("M132").Value = -5015

$ws.Range("H137").Value = 1670.0698
$ws.Range("I137").Value = 1400.7142
$ws.Range("J137").Value = 2172.8667
$ws.Range("K137").Value = 4202.142599999999
$ws.Range("L137").Value = 6518.6001
$ws.Range("M137").Value = -1652.142599999999
$ws.Range("N137").Value = -11618.6001

$ws.Range("H138").Value = 3450.0225
$ws.Range("I138").Value = 1184.3334
$ws.Range("J138").Value = 4436.6934
$ws.Range("K138").Value = 3553.0002
$ws.Range("L138").Value = 13310.0802
$ws.Range("M138").Value = 1586.9998
$ws.Range("N138").Value = -23590.0802

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10000
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10230

$ws.Range("H32").Value = 5015.52
$ws.Range("I32").Value = 3346.725
$ws.Range("J32").Value = 11690.7
$ws.Range("K32").Value = 3346.725
$ws.Range("L32").Value = 11690.7
$ws.Range("M32").Value = -3059.725
$ws.Range("N32").Value = -12264.7

$ws.Range("H52").Value = 27919.334
$ws.Range("J52").Value = 27919.334
$ws.Range("L52").Value = 27919.334
$ws.Range("N52").Value = -28555.334

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H61").Value = 4336.1143
$ws.Range("I61").Value = 4867.357
$ws.Range("J61").Value = 2211.1428
$ws.Range("K61").Value = 4867.357
$ws.Range("L61").Value = 2211.1428
$ws.Range("M61").Value = -4655.357
$ws.Range("N61").Value = -2635.1428

$ws.Range("H74").Value = 1052.6274
$ws.Range("I74").Value = 867
$ws.Range("J74").Value = 1918.8889
$ws.Range("K74").Value = 867
$ws.Range("L74").Value = 1918.8889
$ws.Range("M74").Value = 7
$ws.Range("N74").Value = -3666.8889

$ws.Range("H77").Value = 1052.6274
$ws.Range("I77").Value = 867
$ws.Range("J77").Value = 1918.8889
$ws.Range("K77").Value = 4335
$ws.Range("L77").Value = 9594.4445
$ws.Range("M77").Value = 33
$ws.Range("N77").Value = -18330.4445

$ws.Range("H107").Value = 47614.25
$ws.Range("J107").Value = 47614.25
$ws.Range("L107").Value = 47614.25
$ws.Range("N107").Value = -55294.25

$ws.Range("H109").Value = 48094.5
$ws.Range("J109").Value = 48094.5
$ws.Range("L109").Value = 48094.5
$ws.Range("N109").Value = -50868.5

$ws.Range("H132").Value = 3809.7354
$ws.Range("I132").Value = 2564.4583
$ws.Range("J132").Value = 6798.4
$ws.Range("K132").Value = 7693.374899999999
$ws.Range("L132").Value = 20395.2
$ws.Range("M132").Value = -5163.374899999999
$ws.Range("N132").Value = -25455.2

$ws.Range("H136").Value = 4336.1143
$ws.Range("I136").Value = 4867.357
$ws.Range("J136").Value = 2211.1428
$ws.Range("K136").Value = 14602.071
$ws.Range("L136").Value = 6633.428400000001
$ws.Range("M136").Value = -12052.071
$ws.Range("N136").Value = -11733.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 34483636
$ws.Range("I107").Value = 50000740
$ws.Range("K107").Value = 50000740
$ws.Range("M107").Value = -49998820

$ws.Range("H134").Value = 5455.3438
$ws.Range("I134").Value = 6693.2856
$ws.Range("J134").Value = 3092
$ws.Range("K134").Value = 20079.8568
$ws.Range("L134").Value = 9276
$ws.Range("M134").Value = -17544.8568
$ws.Range("N134").Value = -14346

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 223653.77
$ws.Range("I31").Value = 1409.4822
$ws.Range("J31").Value = 1112630.9
$ws.Range("K31").Value = 1409.4822
$ws.Range("L31").Value = 1112630.9
$ws.Range("M31").Value = -1114.4822
$ws.Range("N31").Value = -1113220.9

$ws.Range("H34").Value = 223653.77
$ws.Range("I34").Value = 1409.4822
$ws.Range("J34").Value = 1112630.9
$ws.Range("K34").Value = 1409.4822
$ws.Range("L34").Value = 1112630.9
$ws.Range("M34").Value = -1207.4822
$ws.Range("N34").Value = -1113034.9

$ws.Range("H99").Value = 5960849.5
$ws.Range("I99").Value = 9324
$ws.Range("J99").Value = 17863900
$ws.Range("K99").Value = 9324
$ws.Range("L99").Value = 17863900
$ws.Range("M99").Value = -7826
$ws.Range("N99").Value = -17866896

$ws.Range("H126").Value = 5960849.5
$ws.Range("I126").Value = 9324
$ws.Range("J126").Value = 17863900
$ws.Range("K126").Value = 27972
$ws.Range("L126").Value = 53591700
$ws.Range("M126").Value = -25502
$ws.Range("N126").Value = -53596640

$ws.Range("H132").Value = 2413.0466
$ws.Range("I132").Value = 1944.4062
$ws.Range("J132").Value = 3776.3635
$ws.Range("K132").Value = 5833.2186
$ws.Range("L132").Value = 11329.0905
$ws.Range("M132").Value = -3303.2186
$ws.Range("N132").Value = -16389.0905

$ws.Range("H134").Value = 2551.675
$ws.Range("I134").Value = 3526.8635
$ws.Range("K134").Value = 10580.5905
$ws.Range("M134").Value = -8045.5905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1898.3846
$ws.Range("J34").Value = 2499.875
$ws.Range("L34").Value = 7499.625
$ws.Range("N34").Value = -7667.625

$ws.Range("H39").Value = 2766.5
$ws.Range("J39").Value = 2766.5
$ws.Range("L39").Value = 8299.5
$ws.Range("N39").Value = -8887.5

$ws.Range("H55").Value = 2275.1155
$ws.Range("J55").Value = 2275.1155
$ws.Range("L55").Value = 6825.3465
$ws.Range("N55").Value = -7179.3465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2375.3257
$ws.Range("I132").Value = 1926.9259
$ws.Range("J132").Value = 3132
$ws.Range("K132").Value = 5780.7777
$ws.Range("L132").Value = 9396
$ws.Range("M132").Value = -3250.7777
$ws.Range("N132").Value = -14456

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 22223098
$ws.Range("I46").Value = 41667256
$ws.Range("J46").Value = 1207.1428
$ws.Range("K46").Value = 41667256
$ws.Range("L46").Value = 1207.1428
$ws.Range("M46").Value = -41667068
$ws.Range("N46").Value = -1583.1428

$ws.Range("H132").Value = 9718917
$ws.Range("I132").Value = 11135555
$ws.Range("J132").Value = 4829
$ws.Range("K132").Value = 33406665
$ws.Range("L132").Value = 14487
$ws.Range("M132").Value = -33404135
$ws.Range("N132").Value = -19547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("K13").Value = 1000
$ws.Range("M13").Value = -860

$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H122").Value = 2590.0322
$ws.Range("I122").Value = 2547.4348
$ws.Range("J122").Value = 2712.5
$ws.Range("K122").Value = 7642.3044
$ws.Range("L122").Value = 8137.5
$ws.Range("M122").Value = -5192.3044
$ws.Range("N122").Value = -13037.5
